# Publish documentation 0.1.1 / ror 0.1.1
# - bump Version and Date on the Metadata sheet
# - add a new "Context" row documenting the element:ContactPoint context,
#   inserted above the existing extension-context row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump version number (row 3: Property="Version")
$ws.Range("B3").Value = "0.1.1"

# Bump publish date (row 8: Property="Date")
$ws.Range("B8").Value = "2023-06-02T12:02:38+02:00"

# The existing row 20 held the only "Context" entry (the extension context).
# Remember its value before touching anything.
$oldContextValue = $ws.Range("B20").Text

# Clone row 20's formatting down into the new row 21 so the appended row
# matches the look of the existing context rows.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)

# Row 20 becomes the new "element:ContactPoint" context …
$ws.Range("A20").Value = "Context"
$ws.Range("B20").Value = "element:ContactPoint"

# … and row 21 keeps the original extension context that used to live in row 20.
$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = $oldContextValue
